$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 327.66666
$ws.Range("I4").Value = 327.66666
$ws.Range("K4").Value = 327.66666
$ws.Range("M4").Value = -213.66666

$ws.Range("H5").Value = 142
$ws.Range("J5").Value = 140
$ws.Range("L5").Value = 140
$ws.Range("N5").Value = -370

$ws.Range("H21").Value = 42333
$ws.Range("I21").Value = 42333
$ws.Range("K21").Value = 42333
$ws.Range("M21").Value = -41865

$ws.Range("H23").Value = 42333
$ws.Range("I23").Value = 42333
$ws.Range("K23").Value = 42333
$ws.Range("M23").Value = -42099

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H82").Value = 2200
$ws.Range("I82").Value = 2200
$ws.Range("K82").Value = 6600
$ws.Range("M82").Value = -6194

$ws.Range("H85").Value = 2200
$ws.Range("I85").Value = 2200
$ws.Range("K85").Value = 6600
$ws.Range("M85").Value = -5196

$ws.Range("H100").Value = 2071
$ws.Range("I100").Value = 1491
$ws.Range("J100").Value = 2361
$ws.Range("K100").Value = 1491
$ws.Range("L100").Value = 2361
$ws.Range("M100").Value = -950
$ws.Range("N100").Value = -3443

$ws.Range("H138").Value = 1698.8
$ws.Range("I138").Value = 997
$ws.Range("J138").Value = 2166.6667
$ws.Range("K138").Value = 2991
$ws.Range("L138").Value = 6500.000100000001
$ws.Range("M138").Value = 2149
$ws.Range("N138").Value = -16780.0001

$ws.Range("H141").Value = 28747.25
$ws.Range("I141").Value = 99989
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 299967
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -294787
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 584.375
$ws.Range("I2").Value = 382.2857
$ws.Range("K2").Value = 382.2857
$ws.Range("M2").Value = -269.2857

$ws.Range("H12").Value = 1400
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1400
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = 1400
$ws.Range("N12").Value = -1746

$ws.Range("H32").Value = 12429.083
$ws.Range("I32").Value = 8794.333000000001
$ws.Range("K32").Value = 8794.333000000001
$ws.Range("M32").Value = -8507.333000000001

$ws.Range("H97").Value = 1235
$ws.Range("I97").Value = 293.75
$ws.Range("K97").Value = 293.75
$ws.Range("M97").Value = 202.25

$ws.Range("H116").Value = 584.375
$ws.Range("I116").Value = 382.2857
$ws.Range("K116").Value = 382.2857
$ws.Range("M116").Value = 1911.7143

$ws.Range("H130").Value = 35449
$ws.Range("J130").Value = 35449
$ws.Range("L130").Value = 35449
$ws.Range("N130").Value = -45489

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 584.375
$ws.Range("I3").Value = 382.2857
$ws.Range("K3").Value = 382.2857
$ws.Range("M3").Value = -268.2857

$ws.Range("H7").Value = 40500
$ws.Range("I7").Value = 66834
$ws.Range("K7").Value = 66834
$ws.Range("M7").Value = -66721

$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -232

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0

$ws.Range("H99").Value = 792
$ws.Range("I99").Value = 910.5
$ws.Range("K99").Value = 910.5
$ws.Range("M99").Value = 587.5

$ws.Range("H134").Value = 2836.8333
$ws.Range("I134").Value = 1004.2
$ws.Range("K134").Value = 3012.6
$ws.Range("M134").Value = -477.6000000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 34.555557
$ws.Range("I7").Value = 19.75
$ws.Range("K7").Value = 19.75
$ws.Range("M7").Value = 93.25

$ws.Range("H28").Value = 33268.418
$ws.Range("I28").Value = 9998
$ws.Range("J28").Value = 35383.91
$ws.Range("K28").Value = 9998
$ws.Range("L28").Value = 35383.91
$ws.Range("M28").Value = -9753
$ws.Range("N28").Value = -35873.91

$ws.Range("H122").Value = 1525
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 11499.5
$ws.Range("I132").Value = 9333
$ws.Range("K132").Value = 27999
$ws.Range("M132").Value = -25469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 25
$ws.Range("I6").Value = 21.428572
$ws.Range("K6").Value = 64.28571599999999
$ws.Range("M6").Value = 48.71428400000001

$ws.Range("H10").Value = 46.125
$ws.Range("J10").Value = 47
$ws.Range("L10").Value = 141
$ws.Range("N10").Value = -419

$ws.Range("H81").Value = 700
$ws.Range("J81").Value = 700
$ws.Range("L81").Value = 2100
$ws.Range("N81").Value = -4346

$ws.Range("H84").Value = 700
$ws.Range("J84").Value = 700
$ws.Range("L84").Value = 6300
$ws.Range("N84").Value = -17532

$ws.Range("H125").Value = 13333.333

$ws.Range("H139").Value = 4620.857
$ws.Range("I139").Value = 2052.1667
$ws.Range("K139").Value = 6156.500100000001
$ws.Range("M139").Value = -1016.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 31
$ws.Range("I2").Value = 29.2
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 29.2
$ws.Range("L2").Value = 34
$ws.Range("M2").Value = 83.8
$ws.Range("N2").Value = -260

$ws.Range("H7").Value = 4000000
$ws.Range("J7").Value = 4000000
$ws.Range("L7").Value = 4000000
$ws.Range("N7").Value = -4000224

$ws.Range("H8").Value = 4000000
$ws.Range("J8").Value = 4000000
$ws.Range("L8").Value = 4000000
$ws.Range("N8").Value = -4000278

$ws.Range("H23").Value = 2013.875
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 500
$ws.Range("M23").Value = -277

$ws.Range("H51").Value = 373999
$ws.Range("J51").Value = 373999
$ws.Range("L51").Value = 373999
$ws.Range("N51").Value = -375017

$ws.Range("H126").Value = 5012
$ws.Range("I126").Value = 5012
$ws.Range("K126").Value = 15036
$ws.Range("M126").Value = -12566

$ws.Range("H132").Value = 2374.5
$ws.Range("I132").Value = 1832.6666
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5497.9998
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2967.9998
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 439.5
$ws.Range("I22").Value = 367
$ws.Range("K22").Value = 367
$ws.Range("M22").Value = -72

$ws.Range("H27").Value = 439.5
$ws.Range("I27").Value = 367
$ws.Range("K27").Value = 367
$ws.Range("M27").Value = -260

$ws.Range("H40").Value = 900
$ws.Range("I40").Value = 900
$ws.Range("K40").Value = 900
$ws.Range("M40").Value = -764

$ws.Range("H82").Value = 573.25
$ws.Range("I82").Value = 500
$ws.Range("K82").Value = 500
$ws.Range("M82").Value = -139

$ws.Range("H85").Value = 573.25
$ws.Range("I85").Value = 500
$ws.Range("K85").Value = 500
$ws.Range("M85").Value = 748

$ws.Range("H128").Value = 35000
$ws.Range("J128").Value = 35000
$ws.Range("L128").Value = 35000
$ws.Range("N128").Value = -44960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("N22").Value = 0

$ws.Range("H96").Value = 2104.8333
$ws.Range("I96").Value = 1535
$ws.Range("K96").Value = 1535
$ws.Range("M96").Value = -162

$ws.Range("H122").Value = 33916.332
$ws.Range("I122").Value = 1750
$ws.Range("K122").Value = 5250
$ws.Range("M122").Value = -2800

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0
